$d = $word.ActiveDocument

# --- Edit 1: merge the four runs in the "Grâce à ma technique..." paragraph ---
$old1 = "Grâce à ma technique de sécurisation des communications, l’information envoyée depuis un émetteur A sera compréhensible au récepteur légitime B mais restera brouillée partout ailleurs."
$new1 = "Grâce à ma technique de sécurisation des communications, l’information envoyée depuis un émetteur A sera compréhensible au récepteur légitime B mais restera brouillée partout ailleurs."

$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

# --- Edit 2: merge the five runs in the "Lorsqu'on émet des signaux..." sentence ---
$old2 = ". Lorsqu’on émet des signaux d’un point A, chacun arrivera en B de manière aléatoire et à des instants différents, en raison des différents obstacles qu’ils rencontrent sur leurs trajectoires (les immeubles, les véhicules,…). Cet aléa est ensuite utilisé pour générer, uniquement entre A et B, une signature spécifique sécurisant la communication. Si un récepteur illégitime est présent à un endroit C quelconque, il réceptionnera l’information provenant de A, "
$new2 = $old2

$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)
